$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Echoes of the Cosmos: Unraveling the Mystery of Dark Energy" "The Rhythm of Chemistry: Unveiling the Symphony of Matter"

# Author name
Replace-Text "Stella Montgomery" "Stephen Coleman"

# Email - first run
Replace-Text "stella" "scoleman@scholar"

# Body paragraph 1
Replace-Text "1) As we peer into the vastness of the universe, searching for answers to its many enigmas, one profound mystery that haunts scientists and astrophysicists alike is the enigmatic substance known as Dark Energy" "The enthralling realm of Chemistry unravels the symphony of matter, revealing the magic within the molecular dance of atoms and elements"
Replace-Text " This mysterious force, believed to be the dominant component of the universe, exerts a profound influence on its fate and evolution" " Like a celestial choreographer, the periodic table organizes a myriad of elements, each possessing its own unique traits and behaviors"
Replace-Text " In this essay, we delve into the depths of Dark Energy, exploring its compelling nature, seeking to understand its role in shaping the destiny of the universe" " Embarking on this captivating journey into the world of Chemistry is like unraveling a symphony in motion, where every compound tells a story of creation and transformation"

# Body paragraph 2
Replace-Text "2) The existence of Dark Energy was first hypothesized to account for the observed acceleration in the expansion of the universe, a phenomenon discovered through meticulous observations of distant galaxies" "Our senses, keen explorers in the chemical realm, delight in the vibrant colors that herald chemical reactions, the tantalizing aromas that waft from molecular interactions, and the tactile sensations that speak of textures and states of matter"
Replace-Text " It is believed to constitute roughly 70% of the total energy density in the universe, dwarfing the contributions of ordinary matter and radiation" " Chemistry's narrative is etched within the pages of history, from the ancient alchemists' quest for the philosopher's stone to the groundbreaking discoveries that have shaped our modern world"
Replace-Text " Despite its profound impact, we have yet to fully comprehend the true nature of Dark Energy, its origin, or its implications for the ultimate destiny of the cosmos" " Chemistry is omnipresent, touching every aspect of our lives, like an intricate tapestry woven into the fabric of our existence"

# Body paragraph 3
Replace-Text "3) One of the most bewildering aspects of Dark Energy is its negative pressure, which leads to its repulsive gravitational force" "In the molecular dance, atoms gracefully waltz, exchanging electrons in a harmonious ballet called chemical bonding"
Replace-Text " The interplay between this repulsive force and the attractive force of ordinary matter gives rise to a unique set of consequences" " These bonds, like invisible threads, forge molecules, the building blocks of our universe"
Replace-Text " The expansion of the universe is accelerated, driving galaxies apart and shaping the large-scale structure of the cosmos" " With each rearrangement of atoms, Chemistry unveils a new creation, whether it be a life-sustaining protein, a resilient plastic, or a dazzling array of pigments that paint the canvas of nature"
Replace-Text " The fate of the universe ultimately rests on the ultimate nature of Dark Energy and whether its influence will cause the expansion to continue indefinitely, leading to a `"Big Rip,`" or whether it will eventually fade away, ushering in an era of eternal cosmic cooling" " From the sparkling snowflakes that adorn winter landscapes to the burst of flavors that tantalize our taste buds, Chemistry's symphony enchants us with its endless creativity"

# Summary
Replace-Text "Comprehending Dark Energy presents one of the most captivating and elusive challenges in contemporary physics" "Chemistry unveils the symphony of matter, revealing the intricate dance of atoms and elements"
Replace-Text " As we continue to unravel the mysteries of this enigmatic substance, we gain a deeper understanding of the universe's grand tapestry" " It enchants us with vibrant colors, tantalizing aromas, and tactile sensations, etching its narrative into the pages of history"
Replace-Text " The quest to " " From the quest for the philosopher's stone to modern-day discoveries, Chemistry's impact is omnipresent, touching every aspect of our "

Write-Output "done basic replacements"
